$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add the "Hyperlink" character style (based on the default paragraph
#    font "a0"), carrying the blue/underline look used by the new link.
# ---------------------------------------------------------------------------
$hlStyle = $d.Styles.Add("Hyperlink", 2)
$hlStyle.BaseStyle = "a0"
$hlStyle.Priority = 99
$hlStyle.UnhideWhenUsed = $true
$hlStyle.Font.Color = 8812614     # RGB(70,120,134) == hex 467886
$hlStyle.Font.Underline = 1       # single underline

# ---------------------------------------------------------------------------
# 2) Mint a real "hyperlink" (external) relationship by adding a throw-away
#    hyperlink at the very end of the document, then removing the scratch
#    text again. The relationship entry itself survives the deletion, and
#    this is the cleanest way to get Word to register the relationship
#    (rather than trying to hand-author package relationship XML).
# ---------------------------------------------------------------------------
$url = "https://github.com/betgws/software"

$scratch = $d.Content
$scratch.Collapse(0)
$scratch.InsertParagraphAfter()
$scratchPara = $d.Paragraphs.Last
$null = $d.Hyperlinks.Add($scratchPara.Range, $url, "", "", $url)
$scratchPara2 = $d.Paragraphs.Last
$scratchPara2.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Paragraph 1 currently holds just a manual line break (<w:br/>). Replace
#    that with the "개인코드: 35246 " text, keeping the paragraph's own
#    (bold) paragraph mark formatting untouched.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Text = "개인코드: 35246 "

# ---------------------------------------------------------------------------
# 4) Insert a brand-new paragraph right after it that holds the hyperlink
#    (display text + a manual line break inside the hyperlink run), matching
#    the bold paragraph-mark formatting of the paragraph it replaces/extends.
# ---------------------------------------------------------------------------
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item(2)

$hyperlinkXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:pPr><w:rPr><w:b w:val="1"/><w:bCs w:val="1"/></w:rPr></w:pPr>
<w:hyperlink r:id="rId9">
<w:r>
<w:rPr>
<w:rFonts w:ascii="맑은 고딕" w:hAnsi="맑은 고딕" w:eastAsia="맑은 고딕" w:cs="맑은 고딕"/>
<w:noProof w:val="0"/>
<w:sz w:val="22"/>
<w:szCs w:val="22"/>
<w:lang w:eastAsia="ko-KR"/>
</w:rPr>
<w:t>https://github.com/betgws/software</w:t>
</w:r>
<w:r><w:br/></w:r>
</w:hyperlink>
</w:p>
'@
$p2.Range.InsertXML($hyperlinkXml)

# ---------------------------------------------------------------------------
# 5) Apply the "Hyperlink" character style to the display-text run. Doing
#    this as a dedicated style assignment (rather than baking <w:rStyle>
#    straight into the inserted XML) is what makes it stick.
# ---------------------------------------------------------------------------
$p2b = $d.Paragraphs.Item(2)
$urlLen = $url.Length
$linkTextRange = $d.Range($p2b.Range.Start, $p2b.Range.Start + $urlLen)
$linkTextRange.Style = "Hyperlink"

Write-Output "done"
